$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain numeric-looking text in the source data
# (e.g. '217.56', thousand-dot formats like '90.996.58', tiny decimals like
# '0.0000249'). Force Text number format before assigning so Excel doesn't
# auto-convert them to numeric cells, matching the original inlineStr text cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.996.58'
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.167.58'
$ws.Range("E3").Value = '  +2.96%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.56'
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '628.65'
$ws.Range("E6").Value = '  +2.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.17'
$ws.Range("E7").Value = '  +30.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.369'
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.164.04'
$ws.Range("E10").Value = '  +2.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.769'
$ws.Range("E11").Value = '  +13.39%  '
$ws.Range("E12").Value = '  +8.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  +4.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.69'
$ws.Range("E14").Value = '  +5.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '35.38'
$ws.Range("E15").Value = '  +9.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.813.38'
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.753.09'
$ws.Range("E17").Value = '  +3.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.176.98'
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.69'
$ws.Range("E19").Value = '  +9.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000218'
$ws.Range("E20").Value = '  +2.92%  '
$ws.Range("E21").Value = '  +6.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '458.48'
$ws.Range("E22").Value = '  +6.74%  '
$ws.Range("E23").Value = '  +9.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.20'
$ws.Range("E24").Value = '  +3.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.99'
$ws.Range("E25").Value = '  +9.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.31'
$ws.Range("E26").Value = '  +7.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.20'
$ws.Range("E27").Value = '  +2.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.340.30'
$ws.Range("E28").Value = '  +3.11%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.33'
$ws.Range("E30").Value = '  +13.13%  '
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("E32").Value = '  -7.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '527.40'
$ws.Range("E33").Value = '  +4.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.97'
$ws.Range("E34").Value = '  +14.10%  '
$ws.Range("E35").Value = '  +32.50%  '
$ws.Range("E36").Value = '  +2.27%  '
$ws.Range("E37").Value = '  +9.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.92'
$ws.Range("E38").Value = '  +7.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.95'
$ws.Range("E39").Value = '  +4.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("E40").Value = '  +5.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0918'
$ws.Range("E41").Value = '  +32.02%  '
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.429'
$ws.Range("E42").Value = '  +16.19%  '
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.23'
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.96'
$ws.Range("E45").Value = '  +6.64%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '148.11'
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.36'
$ws.Range("E48").Value = '  +12.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.49'
$ws.Range("E49").Value = '  +9.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.66'
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.656'
$ws.Range("E51").Value = '  +11.60%  '
